$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -13
$ws.Range("F3").Value = -5
$ws.Range("F5").Value = 2
$ws.Range("F7").Value = -4
